# CASOS_DE_PRUEBA.xlsx
# Commit: "Se agrego mensajes-guia y validar multiples apellidos"
#
# The sheet "Hoja1" has two Excel Tables:
#   Tabla1  A1:D10  -> APELLIDO / NOMBRE / NOTA 1 / NOTA 2
#   Tabla2  G1:G10  -> RESULTADOS ESPERADOS
# Row 10 held a "sentinel" row (A10 = 0, G10 = "FIN DEL PROGRAMA") that
# marked the end of the data used by the lookup formulas on the other
# (hidden/not-present-here) sheet. This change inserts one more test case
# -- a student with a compound/multi-word surname ("De la Cruz") -- right
# above that sentinel row, so both tables grow by one row (now ending on
# row 11), and the sentinel moves down to row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 10, pushing the old row 10 (the sentinel) down to
# row 11. Then strip any formatting Excel may have copied down onto the
# new row so the new cells start out with the default style.
$ws.Rows(10).Insert()
$ws.Range("A10:D10").ClearFormats()

# New test case: a student whose APELLIDO is two words, to validate the
# lookup formulas handle multi-word surnames correctly.
$ws.Range("A10").Value = "De la Cruz"
$ws.Range("B10").Value = "Esteban"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 3
$ws.Range("G10").Value = "REPRUEBA"

# Grow both tables so the new row (and the resident sentinel row, now at
# 11) stay inside the table ranges.
$lo1 = $ws.ListObjects.Item("Tabla1")
$lo2 = $ws.ListObjects.Item("Tabla2")
$lo1.Resize($ws.Range("A1:D11"))
$lo2.Resize($ws.Range("G1:G11"))

# "Martinez" (A9) had a leftover underline style; drop it as part of the
# formatting cleanup in this commit.
$ws.Range("A9").Font.Underline = $False

# Leave the selection where the author left it after the edit.
$ws.Range("A13").Select()
